$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(20, 8).Value = 20000
$ws.Cells.Item(20, 9).Value = 20000
$ws.Cells.Item(20, 11).Value = 20000
$ws.Cells.Item(20, 13).Value = -19770

$ws.Cells.Item(26, 8).Value = 19980
$ws.Cells.Item(26, 10).Value = 21250
$ws.Cells.Item(26, 12).Value = 21250
$ws.Cells.Item(26, 14).Value = -21938

$ws.Cells.Item(35, 8).Value = 20000
$ws.Cells.Item(35, 9).Value = 20000
$ws.Cells.Item(35, 11).Value = 20000
$ws.Cells.Item(35, 13).Value = -19621

$ws.Cells.Item(41, 8).Value = 425.14285
$ws.Cells.Item(41, 9).Value = 371
$ws.Cells.Item(41, 10).Value = 750
$ws.Cells.Item(41, 11).Value = 371
$ws.Cells.Item(41, 12).Value = 750
$ws.Cells.Item(41, 13).Value = 69
$ws.Cells.Item(41, 14).Value = -1630

$ws.Cells.Item(58, 8).Value = 365.08334

$ws.Cells.Item(86, 8).Value = 13918
$ws.Cells.Item(86, 9).Value = 26172.5
$ws.Cells.Item(86, 11).Value = 26172.5
$ws.Cells.Item(86, 13).Value = -25049.5

$ws.Cells.Item(88, 8).Value = 15213215
$ws.Cells.Item(88, 10).Value = 16595962
$ws.Cells.Item(88, 12).Value = 16595962
$ws.Cells.Item(88, 14).Value = -16596774

$ws.Cells.Item(89, 8).Value = 13918
$ws.Cells.Item(89, 9).Value = 26172.5
$ws.Cells.Item(89, 11).Value = 130862.5
$ws.Cells.Item(89, 13).Value = -125246.5

$ws.Cells.Item(91, 8).Value = 15213215
$ws.Cells.Item(91, 10).Value = 16595962
$ws.Cells.Item(91, 12).Value = 16595962
$ws.Cells.Item(91, 14).Value = -16598770

$ws.Cells.Item(98, 8).Value = 56819336
$ws.Cells.Item(98, 9).Value = 89286690
$ws.Cells.Item(98, 10).Value = 1469.75
$ws.Cells.Item(98, 11).Value = 89286690
$ws.Cells.Item(98, 12).Value = 1469.75
$ws.Cells.Item(98, 13).Value = -89285192
$ws.Cells.Item(98, 14).Value = -4465.75

$ws.Cells.Item(103, 8).Value = 597.9231
$ws.Cells.Item(103, 9).Value = 663.25
$ws.Cells.Item(103, 10).Value = 568.8889
$ws.Cells.Item(103, 11).Value = 1989.75
$ws.Cells.Item(103, 12).Value = 1706.6667
$ws.Cells.Item(103, 13).Value = -1403.75
$ws.Cells.Item(103, 14).Value = -2878.6667

$ws.Cells.Item(106, 8).Value = 83335576
$ws.Cells.Item(106, 9).Value = 166667730
$ws.Cells.Item(106, 11).Value = 166667730
$ws.Cells.Item(106, 13).Value = -166667099

$ws.Cells.Item(122, 8).Value = 56819336
$ws.Cells.Item(122, 9).Value = 89286690
$ws.Cells.Item(122, 10).Value = 1469.75
$ws.Cells.Item(122, 11).Value = 267860070
$ws.Cells.Item(122, 12).Value = 4409.25
$ws.Cells.Item(122, 13).Value = -267857620
$ws.Cells.Item(122, 14).Value = -9309.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 678.0909
$ws.Cells.Item(2, 9).Value = 557.375
$ws.Cells.Item(2, 11).Value = 557.375
$ws.Cells.Item(2, 13).Value = -444.375

$ws.Cells.Item(32, 8).Value = 14094410
$ws.Cells.Item(32, 9).Value = 17865324
$ws.Cells.Item(32, 10).Value = 16332
$ws.Cells.Item(32, 11).Value = 17865324
$ws.Cells.Item(32, 12).Value = 16332
$ws.Cells.Item(32, 13).Value = -17865037
$ws.Cells.Item(32, 14).Value = -16906

$ws.Cells.Item(45, 8).Value = 48705.668
$ws.Cells.Item(45, 9).Value = 84159.914
$ws.Cells.Item(45, 11).Value = 84159.914
$ws.Cells.Item(45, 13).Value = -83782.914

$ws.Cells.Item(110, 8).Value = 1319.3
$ws.Cells.Item(110, 9).Value = 783.087
$ws.Cells.Item(110, 11).Value = 783.087
$ws.Cells.Item(110, 13).Value = 1261.913

$ws.Cells.Item(116, 8).Value = 678.0909
$ws.Cells.Item(116, 9).Value = 557.375
$ws.Cells.Item(116, 11).Value = 557.375
$ws.Cells.Item(116, 13).Value = 1736.625

$ws.Cells.Item(132, 8).Value = 27782076
$ws.Cells.Item(132, 10).Value = 12311.1
$ws.Cells.Item(132, 12).Value = 36933.3
$ws.Cells.Item(132, 14).Value = -41993.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 678.0909
$ws.Cells.Item(3, 9).Value = 557.375
$ws.Cells.Item(3, 11).Value = 557.375
$ws.Cells.Item(3, 13).Value = -443.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 5383.6665
$ws.Cells.Item(6, 9).Value = 6060
$ws.Cells.Item(6, 10).Value = 2002
$ws.Cells.Item(6, 11).Value = 6060
$ws.Cells.Item(6, 12).Value = 2002
$ws.Cells.Item(6, 13).Value = -5947
$ws.Cells.Item(6, 14).Value = -2228

$ws.Cells.Item(16, 8).Value = 1399.5853
$ws.Cells.Item(16, 9).Value = 721.5833
$ws.Cells.Item(16, 10).Value = 2356.7646
$ws.Cells.Item(16, 11).Value = 721.5833
$ws.Cells.Item(16, 12).Value = 2356.7646
$ws.Cells.Item(16, 13).Value = -434.5833
$ws.Cells.Item(16, 14).Value = -2930.7646

$ws.Cells.Item(99, 8).Value = 71438450
$ws.Cells.Item(99, 9).Value = 125014616
$ws.Cells.Item(99, 10).Value = 3548.3333
$ws.Cells.Item(99, 11).Value = 125014616
$ws.Cells.Item(99, 12).Value = 3548.3333
$ws.Cells.Item(99, 13).Value = -125013118
$ws.Cells.Item(99, 14).Value = -6544.3333

$ws.Cells.Item(105, 8).Value = 5474.5654
$ws.Cells.Item(105, 9).Value = 6218.6113
$ws.Cells.Item(105, 10).Value = 2796
$ws.Cells.Item(105, 11).Value = 6218.6113
$ws.Cells.Item(105, 12).Value = 2796
$ws.Cells.Item(105, 13).Value = -4471.6113
$ws.Cells.Item(105, 14).Value = -6290

$ws.Cells.Item(113, 8).Value = 1399.5853
$ws.Cells.Item(113, 9).Value = 721.5833
$ws.Cells.Item(113, 10).Value = 2356.7646
$ws.Cells.Item(113, 11).Value = 721.5833
$ws.Cells.Item(113, 12).Value = 2356.7646
$ws.Cells.Item(113, 13).Value = 1448.4167
$ws.Cells.Item(113, 14).Value = -6696.7646

$ws.Cells.Item(126, 8).Value = 71438450
$ws.Cells.Item(126, 9).Value = 125014616
$ws.Cells.Item(126, 10).Value = 3548.3333
$ws.Cells.Item(126, 11).Value = 375043848
$ws.Cells.Item(126, 12).Value = 10644.9999
$ws.Cells.Item(126, 13).Value = -375041378
$ws.Cells.Item(126, 14).Value = -15584.9999

$ws.Cells.Item(132, 8).Value = 13894912
$ws.Cells.Item(132, 9).Value = 1316.4375
$ws.Cells.Item(132, 10).Value = 41682104
$ws.Cells.Item(132, 11).Value = 3949.3125
$ws.Cells.Item(132, 12).Value = 125046312
$ws.Cells.Item(132, 13).Value = -1419.3125
$ws.Cells.Item(132, 14).Value = -125051372

$ws.Cells.Item(134, 8).Value = 1473.4117
$ws.Cells.Item(134, 9).Value = 1192.25
$ws.Cells.Item(134, 10).Value = 1723.3334
$ws.Cells.Item(134, 11).Value = 3576.75
$ws.Cells.Item(134, 12).Value = 5170.0002
$ws.Cells.Item(134, 13).Value = -1041.75
$ws.Cells.Item(134, 14).Value = -10240.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1172.8182
$ws.Cells.Item(34, 10).Value = 1999.8
$ws.Cells.Item(34, 12).Value = 5999.4
$ws.Cells.Item(34, 14).Value = -6167.4

$ws.Cells.Item(39, 8).Value = 1209.4546
$ws.Cells.Item(39, 9).Value = 387.5
$ws.Cells.Item(39, 10).Value = 1679.1428
$ws.Cells.Item(39, 11).Value = 1162.5
$ws.Cells.Item(39, 12).Value = 5037.428400000001
$ws.Cells.Item(39, 13).Value = -868.5
$ws.Cells.Item(39, 14).Value = -5625.428400000001

$ws.Cells.Item(55, 8).Value = 1650
$ws.Cells.Item(55, 10).Value = 1775
$ws.Cells.Item(55, 12).Value = 5325
$ws.Cells.Item(55, 14).Value = -5679

$ws.Cells.Item(122, 8).Value = 16670594
$ws.Cells.Item(122, 9).Value = 125000390
$ws.Cells.Item(122, 10).Value = 4471.5
$ws.Cells.Item(122, 11).Value = 1125003510
$ws.Cells.Item(122, 12).Value = 40243.5
$ws.Cells.Item(122, 13).Value = -1125001060
$ws.Cells.Item(122, 14).Value = -45143.5

$ws.Cells.Item(131, 8).Value = 696.4400000000001
$ws.Cells.Item(131, 10).Value = 774.675
$ws.Cells.Item(131, 12).Value = 2324.025
$ws.Cells.Item(131, 14).Value = -12404.025

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 10002260
$ws.Cells.Item(80, 9).Value = 2512.875
$ws.Cells.Item(80, 11).Value = 2512.875
$ws.Cells.Item(80, 13).Value = -1514.875

$ws.Cells.Item(83, 8).Value = 10002260
$ws.Cells.Item(83, 9).Value = 2512.875
$ws.Cells.Item(83, 11).Value = 12564.375
$ws.Cells.Item(83, 13).Value = -7572.375

$ws.Cells.Item(97, 8).Value = 733.3333
$ws.Cells.Item(97, 9).Value = 713.5
$ws.Cells.Item(97, 11).Value = 713.5
$ws.Cells.Item(97, 13).Value = -217.5

$ws.Cells.Item(122, 8).Value = 15630525
$ws.Cells.Item(122, 9).Value = 20839660
$ws.Cells.Item(122, 11).Value = 62518980
$ws.Cells.Item(122, 13).Value = -62516530

$ws.Cells.Item(132, 8).Value = 19313
$ws.Cells.Item(132, 9).Value = 880
$ws.Cells.Item(132, 11).Value = 2640
$ws.Cells.Item(132, 13).Value = -110

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2562.5
$ws.Cells.Item(7, 9).Value = 2500
$ws.Cells.Item(7, 10).Value = 3000
$ws.Cells.Item(7, 11).Value = 2500
$ws.Cells.Item(7, 12).Value = 3000
$ws.Cells.Item(7, 13).Value = -2388
$ws.Cells.Item(7, 14).Value = -3224

$ws.Cells.Item(46, 8).Value = 3789922.8
$ws.Cells.Item(46, 9).Value = 8333953
$ws.Cells.Item(46, 10).Value = 3230.6667
$ws.Cells.Item(46, 11).Value = 8333953
$ws.Cells.Item(46, 12).Value = 3230.6667
$ws.Cells.Item(46, 13).Value = -8333765
$ws.Cells.Item(46, 14).Value = -3606.6667

$ws.Cells.Item(126, 8).Value = 2562.5
$ws.Cells.Item(126, 9).Value = 2500
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 7500
$ws.Cells.Item(126, 12).Value = 9000
$ws.Cells.Item(126, 13).Value = -5030
$ws.Cells.Item(126, 14).Value = -13940

$ws.Cells.Item(132, 8).Value = 19612948
$ws.Cells.Item(132, 9).Value = 41668840
$ws.Cells.Item(132, 11).Value = 125006520
$ws.Cells.Item(132, 13).Value = -125003990

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 52632188
$ws.Cells.Item(113, 9).Value = 71429070
$ws.Cells.Item(113, 10).Value = 928
$ws.Cells.Item(113, 11).Value = 214287210
$ws.Cells.Item(113, 12).Value = 2784
$ws.Cells.Item(113, 13).Value = -214285040
$ws.Cells.Item(113, 14).Value = -7124

$ws.Cells.Item(126, 8).Value = 1817
$ws.Cells.Item(126, 10).Value = 2508.3333
$ws.Cells.Item(126, 12).Value = 7524.999899999999
$ws.Cells.Item(126, 14).Value = -12464.9999

$ws.Cells.Item(132, 8).Value = 19480.258
$ws.Cells.Item(132, 9).Value = 27573.719
$ws.Cells.Item(132, 11).Value = 82721.15700000001
$ws.Cells.Item(132, 13).Value = -80191.15700000001
